# Auto-generated update of Sargatanas_Profits.xlsx market-price columns
# (H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#  K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1055.6154
$ws.Range("I33").Value = 1291.7778
$ws.Range("K33").Value = 1291.7778
$ws.Range("M33").Value = -1062.7778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2803.8
$ws.Range("I132").Value = 2759.532
$ws.Range("K132").Value = 8278.596000000001
$ws.Range("M132").Value = -5748.596000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3750.0557
$ws.Range("I137").Value = 7200.6665
$ws.Range("K137").Value = 21601.9995
$ws.Range("M137").Value = -19051.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4924.6875
$ws.Range("I138").Value = 976.125
$ws.Range("J138").Value = 8873.25
$ws.Range("K138").Value = 2928.375
$ws.Range("L138").Value = 26619.75
$ws.Range("M138").Value = 2211.625
$ws.Range("N138").Value = -36899.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2275966.8
$ws.Range("I32").Value = 2503108.5
$ws.Range("J32").Value = 4549.4
$ws.Range("K32").Value = 2503108.5
$ws.Range("L32").Value = 4549.4
$ws.Range("M32").Value = -2502821.5
$ws.Range("N32").Value = -5123.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7190.1816
$ws.Range("I61").Value = 3156
$ws.Range("J61").Value = 14250
$ws.Range("K61").Value = 3156
$ws.Range("L61").Value = 14250
$ws.Range("M61").Value = -2944
$ws.Range("N61").Value = -14674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7191.5557
$ws.Range("I132").Value = 5785.269
$ws.Range("K132").Value = 17355.807
$ws.Range("M132").Value = -14825.807

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7190.1816
$ws.Range("I136").Value = 3156
$ws.Range("J136").Value = 14250
$ws.Range("K136").Value = 9468
$ws.Range("L136").Value = 42750
$ws.Range("M136").Value = -6918
$ws.Range("N136").Value = -47850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4001.1528
$ws.Range("I134").Value = 1970.9818
$ws.Range("J134").Value = 10569.353
$ws.Range("K134").Value = 5912.945400000001
$ws.Range("L134").Value = 31708.059
$ws.Range("M134").Value = -3377.945400000001
$ws.Range("N134").Value = -36778.05899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6133.75
$ws.Range("I31").Value = 2638.6191
$ws.Range("J31").Value = 11779.73
$ws.Range("K31").Value = 2638.6191
$ws.Range("L31").Value = 11779.73
$ws.Range("M31").Value = -2343.6191
$ws.Range("N31").Value = -12369.73

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6133.75
$ws.Range("I34").Value = 2638.6191
$ws.Range("J34").Value = 11779.73
$ws.Range("K34").Value = 2638.6191
$ws.Range("L34").Value = 11779.73
$ws.Range("M34").Value = -2436.6191
$ws.Range("N34").Value = -12183.73

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12505817
$ws.Range("I58").Value = 23811226
$ws.Range("K58").Value = 23811226
$ws.Range("M58").Value = -23811023

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4689.067
$ws.Range("I132").Value = 1872.4736
$ws.Range("J132").Value = 9554.091
$ws.Range("K132").Value = 5617.4208
$ws.Range("L132").Value = 28662.273
$ws.Range("M132").Value = -3087.4208
$ws.Range("N132").Value = -33722.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 12505817
$ws.Range("I136").Value = 23811226
$ws.Range("K136").Value = 71433678
$ws.Range("M136").Value = -71431128

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8005000
$ws.Range("I5").Value = 20002000
$ws.Range("J5").Value = 7000
$ws.Range("K5").Value = 60006000
$ws.Range("L5").Value = 21000
$ws.Range("M5").Value = -60005888
$ws.Range("N5").Value = -21224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 3591.3333
$ws.Range("I24").Value = 3490
$ws.Range("J24").Value = 4250
$ws.Range("K24").Value = 10470
$ws.Range("L24").Value = 12750
$ws.Range("M24").Value = -10240
$ws.Range("N24").Value = -13210

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 5500
$ws.Range("I58").Value = 5500
$ws.Range("K58").Value = 16500
$ws.Range("M58").Value = -16372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 8005000
$ws.Range("I135").Value = 20002000
$ws.Range("J135").Value = 7000
$ws.Range("K135").Value = 180018000
$ws.Range("L135").Value = 63000
$ws.Range("M135").Value = -180015465
$ws.Range("N135").Value = -68070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1746.2368
$ws.Range("I97").Value = 1686.7391
$ws.Range("J97").Value = 1837.4667
$ws.Range("K97").Value = 1686.7391
$ws.Range("L97").Value = 1837.4667
$ws.Range("M97").Value = -1190.7391
$ws.Range("N97").Value = -2829.4667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4152.6772
$ws.Range("I102").Value = 4133.8076
$ws.Range("K102").Value = 4133.8076
$ws.Range("M102").Value = -2511.8076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2514.611
$ws.Range("J126").Value = 2528.8
$ws.Range("L126").Value = 7586.400000000001
$ws.Range("N126").Value = -12526.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2176.0334
$ws.Range("I132").Value = 2176.0334
$ws.Range("K132").Value = 6528.100199999999
$ws.Range("M132").Value = -3998.100199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 25709
$ws.Range("I133").Value = 25709
$ws.Range("K133").Value = 25709
$ws.Range("M133").Value = -20649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6143.7856
$ws.Range("I7").Value = 4200.6
$ws.Range("J7").Value = 7223.3335
$ws.Range("K7").Value = 4200.6
$ws.Range("L7").Value = 7223.3335
$ws.Range("M7").Value = -4088.6
$ws.Range("N7").Value = -7447.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13598.556
$ws.Range("I22").Value = 993.6667
$ws.Range("J22").Value = 19901
$ws.Range("K22").Value = 993.6667
$ws.Range("L22").Value = 19901
$ws.Range("M22").Value = -698.6667
$ws.Range("N22").Value = -20491

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 13598.556
$ws.Range("I27").Value = 993.6667
$ws.Range("J27").Value = 19901
$ws.Range("K27").Value = 993.6667
$ws.Range("L27").Value = 19901
$ws.Range("M27").Value = -886.6667
$ws.Range("N27").Value = -20115

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5417.625
$ws.Range("I40").Value = 4717.9
$ws.Range("K40").Value = 4717.9
$ws.Range("M40").Value = -4581.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3978.8845
$ws.Range("I122").Value = 3637.5945
$ws.Range("K122").Value = 10912.7835
$ws.Range("M122").Value = -8462.783500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6143.7856
$ws.Range("I126").Value = 4200.6
$ws.Range("J126").Value = 7223.3335
$ws.Range("K126").Value = 12601.8
$ws.Range("L126").Value = 21670.0005
$ws.Range("M126").Value = -10131.8
$ws.Range("N126").Value = -26610.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13895107
$ws.Range("I132").Value = 20837104
$ws.Range("K132").Value = 62511312
$ws.Range("M132").Value = -62508782

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6168.25
$ws.Range("I122").Value = 3999
$ws.Range("K122").Value = 11997
$ws.Range("M122").Value = -9547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15652939
$ws.Range("I132").Value = 17862252
$ws.Range("K132").Value = 53586756
$ws.Range("M132").Value = -53584226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 136166.67
$ws.Range("J133").Value = 136166.67
$ws.Range("L133").Value = 136166.67
$ws.Range("N133").Value = -146286.67

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 50005012
$ws.Range("I136").Value = 125001380
$ws.Range("J136").Value = 7434
$ws.Range("K136").Value = 375004140
$ws.Range("L136").Value = 22302
$ws.Range("M136").Value = -375001590
$ws.Range("N136").Value = -27402
